# Expansao das analises automaticas:
# - adiciona 3 novas colunas (L: apoio_medio, M: contribuicoes, N: media_contribuicoes)
# - corrige escala das colunas E (particip) e F (taxa_sucesso), de fracao (0-1) para
#   percentual numerico (0-100), mantendo o formato de celula ja existente

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Novos cabecalhos (linha 1), copiando o formato do cabecalho existente (K1) ---
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- Linha 2 ---
$ws.Range("E2").Value = 77.82771535580524
$ws.Range("F2").Value = 62.65640038498557
$ws.Range("L2").Value = 92.37335010030465
$ws.Range("M2").Value = 202024
$ws.Range("N2").Value = 310.3287250384025

# --- Linha 3 ---
$ws.Range("E3").Value = 22.17228464419476
$ws.Range("F3").Value = 60.47297297297297
$ws.Range("L3").Value = 87.79024763378376
$ws.Range("M3").Value = 61529
$ws.Range("N3").Value = 343.7374301675978

# --- Linha 4 ---
$ws.Range("E4").Value = 76.02179836512262
$ws.Range("F4").Value = 93.63799283154121
$ws.Range("L4").Value = 89.64725872903666
$ws.Range("M4").Value = 145834
$ws.Range("N4").Value = 139.5540669856459

# --- Linha 5 ---
$ws.Range("E5").Value = 23.97820163487738
$ws.Range("F5").Value = 96.02272727272727
$ws.Range("L5").Value = 91.47778329877478
$ws.Range("M5").Value = 57812
$ws.Range("N5").Value = 171.0414201183432

# --- Linha 6 ---
$ws.Range("E6").Value = 89.61988304093568
$ws.Range("F6").Value = 22.02283849918434
$ws.Range("L6").Value = 19.42752789799076
$ws.Range("M6").Value = 2118
$ws.Range("N6").Value = 15.68888888888889

# --- Linha 7 ---
$ws.Range("E7").Value = 10.38011695906433
$ws.Range("F7").Value = 23.94366197183098
$ws.Range("L7").Value = 22.66059629822632
$ws.Range("M7").Value = 90
$ws.Range("N7").Value = 5.294117647058823

Write-Host "edit applied"
